# Update "想去人数" (interest count) values on the "展览" and "全部类型" sheets
# F2: 616 -> 618
# F3: 3763 -> 3779
# F4: 103 -> 104
# F5: 718 -> 719

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 618
    $ws.Range("F3").Value = 3779
    $ws.Range("F4").Value = 104
    $ws.Range("F5").Value = 719
}
